$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row before the existing row 101. This shifts the
# previous rows 101..148 down to 102..149 (so the old last data row, 148,
# becomes 149, and the old 147 becomes 148 -- matching the target diff).
$ws.Rows.Item(101).Insert()

# Clone the row immediately below (which now holds the data that used to be
# in row 101) into the freshly-inserted blank row 101, so all the columns
# that stay constant across this data set (Mercado ID, Mercado, Region,
# Codreg, Categoria ID, Categoria, Variedad, Calidad, Unidad de
# comercializacion, Origen, Kg o Unidades, Clasificacion) plus the date
# column's style come along for free.
$ws.Rows.Item(102).Copy()
$ws.Rows.Item(101).PasteSpecial()

# Now overwrite the cells that actually hold new data for this new record.
$ws.Range("D101").Value = 44518
$ws.Range("J101").Value = 125
$ws.Range("K101").Value = 5500
$ws.Range("L101").Value = 6000
$ws.Range("M101").Value = 5740
$ws.Range("P101").Value = 359
